# Major update: if a trip is longer than 90 mins, split it into two groups.
# This renumbers the route/group labels in columns A ("label_x") and L
# ("Ruta") for rows 2-22 (routes 1-9) down by one (0-based numbering), and
# corrects a few travel-time values (columns I/K) on rows 23, 31 and 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-22: decrement the route number in column A and column L by 1.
for ($r = 2; $r -le 22; $r++) {
    $current = $ws.Cells.Item($r, 1).Value2
    $newVal = $current - 1
    $ws.Cells.Item($r, 1).Value = $newVal
    $ws.Cells.Item($r, 12).Value = $newVal
}

# Correct travel-time figures (columns I and K) on a few rows.
$ws.Range("I23").Value = 19.2
$ws.Range("K23").Value = 19.2

$ws.Range("I31").Value = 19.2
$ws.Range("K31").Value = 19.2

$ws.Range("I37").Value = 19.9
$ws.Range("K37").Value = 19.9
